# Update "想去人数" (want-to-go count) values per the diff, output generated at 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 3337
$ws.Range("F6").Value = 1117
$ws.Range("F7").Value = 2240
$ws.Range("F8").Value = 2158
$ws.Range("F9").Value = 1129
$ws.Range("F10").Value = 615
$ws.Range("F12").Value = 1690
$ws.Range("F17").Value = 239
$ws.Range("F18").Value = 1604
$ws.Range("F19").Value = 263
$ws.Range("F20").Value = 654
$ws.Range("F22").Value = 263
$ws.Range("F24").Value = 12346
$ws.Range("F25").Value = 12399
$ws.Range("F27").Value = 712
$ws.Range("F29").Value = 252
$ws.Range("F36").Value = 615

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 41
$ws.Range("F8").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 88

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 3337
$ws.Range("F7").Value = 1117
$ws.Range("F8").Value = 2240
$ws.Range("F9").Value = 2158
$ws.Range("F10").Value = 1129
$ws.Range("F11").Value = 615
$ws.Range("F12").Value = 88
$ws.Range("F14").Value = 1690
$ws.Range("F21").Value = 41
$ws.Range("F22").Value = 239
$ws.Range("F23").Value = 1604
$ws.Range("F24").Value = 263
$ws.Range("F25").Value = 654
$ws.Range("F27").Value = 263
$ws.Range("F29").Value = 12346
$ws.Range("F30").Value = 12399
$ws.Range("F32").Value = 712
$ws.Range("F34").Value = 252
$ws.Range("F43").Value = 615
$ws.Range("F45").Value = 4

